# Apply edit: add "Max Preps" column (I) to the "Prof" sheet with values
# for each existing teacher row, as part of including all input variables
# to be passed to matlab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prof")

$ws.Range("I1").Value = "Max Preps"
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 2

# Move the active selection, matching the saved view state of the sheet.
$ws.Activate()
$ws.Range("J4").Select()
